# Applies the "Fixes reference to figure" edit to the single slide of
# ossmodelv2.pptx: re-positions/re-sizes several pictures, connectors and
# numbered ovals, re-points a connector's start endpoint, and bumps the
# numbered ovals' label text to 24pt bold.
#
# NOTE on the literal numbers below: Shape.Left/Top/Width/Height round-trip
# through a 32-bit (Single) float in this object model (same as real
# PowerPoint VBA), so the plain EMU/12700 point value can land one EMU off
# after conversion back on save. The constants here were solved so that,
# after the float32 round-trip, they serialize to the exact target EMU
# values from the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

# --- Pictures: move only (size unchanged) ------------------------------

$shape62 = Get-ShapeByName $s "Shape 62"
$shape62.Left = 474.3304724409449

$shape81 = Get-ShapeByName $s "Shape 81"
$shape81.Left = 288.1592913385827
$shape81.Top = 354.60072336141735

$shape89pic = Get-ShapeByName $s "Shape 89"
$shape89pic.Left = 319.8248818897638
$shape89pic.Top = 438.5581207362205

$shape79 = Get-ShapeByName $s "Shape 79"
$shape79.Left = 213.08984251968505
$shape79.Top = 354.88102722204724

$shape87 = Get-ShapeByName $s "Shape 87"
$shape87.Left = 168.56299592598427
$shape87.Top = 391.14103702204727

$picture52 = Get-ShapeByName $s "Picture 52"
$picture52.Left = 170.6915818031496
$picture52.Top = 428.96583557165354

# --- Connectors ----------------------------------------------------------

# "Straight Arrow Connector 62" (id 63): now starts at Shape 62 (id 14,
# connection site 1); drop the vertical flip and resize/move it.
$conn62 = Get-ShapeByName $s "Straight Arrow Connector 62"
$conn62.ConnectorFormat.BeginConnect($shape62, 1)
$conn62.VerticalFlip = 0
$conn62.Left = 347.19173228346455
$conn62.Top = 401.6790619181102
$conn62.Width = 127.13874015748031
$conn62.Height = 0.49291338582677163

# "Straight Arrow Connector 88" (id 89): now starts at Shape 61 (id 5,
# inside "Group 70", connection site 4) instead of ending at Shape 62.
$group70 = Get-ShapeByName $s "Group 70"
$shape61 = $null
for ($i = 1; $i -le $group70.GroupItems.Count; $i++) {
    $item = $group70.GroupItems.Item($i)
    if ($item.Name -eq "Shape 61") { $shape61 = $item }
}
$conn88 = Get-ShapeByName $s "Straight Arrow Connector 88"
$conn88.ConnectorFormat.EndDisconnect()
$conn88.ConnectorFormat.BeginConnect($shape61, 4)
$conn88.Left = 609.6630708661418
$conn88.Top = 191.61338582677166
$conn88.Width = 196.36890413779528
$conn88.Height = 228.54496062992126

# "Straight Arrow Connector 34" (id 35): resize only.
$conn34 = Get-ShapeByName $s "Straight Arrow Connector 34"
$conn34.Width = 1.2781102362204724
$conn34.Height = 142.98149876299215

# --- Numbered ovals: move/resize + bold 24pt label ----------------------

function Update-Oval($slide, $name, $left, $top, $width, $height) {
    $sh = Get-ShapeByName $slide $name
    $sh.Left = $left
    $sh.Top = $top
    $sh.Width = $width
    $sh.Height = $height
    $tr = $sh.TextFrame.TextRange
    $tr.Font.Size = 24
    $tr.Font.Bold = $true
}

Update-Oval $s "Oval 15" 436.8520472440945 172.8 39.514410048818895 37.598268516535434
Update-Oval $s "Oval 41" 606.3128357456693 288.9971771543307 39.514410048818895 37.598268516535434
Update-Oval $s "Oval 42" 786.7625427850394 44.12527659055118 39.514410048818895 37.598268516535434
Update-Oval $s "Oval 48" 343.6582795165354 354.14535433070864 39.514410048818895 37.598268516535434
